$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D retains text formatting so numeric-looking values
# (e.g. "0.999", "1.20", "418.77") are not coerced into numbers,
# which would lose trailing zeros / introduce float rounding noise.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '68.825.54'
$ws.Range("E2").Value = '  +0.70%  '
$ws.Range("D3").Value = '3.846.89'
$ws.Range("E3").Value = '  -1.53%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '523.41'
$ws.Range("E5").Value = '  +7.34%  '
$ws.Range("D6").Value = '141.90'
$ws.Range("E6").Value = '  -2.68%  '
$ws.Range("D7").Value = '0.605'
$ws.Range("E7").Value = '  -2.84%  '
$ws.Range("D8").Value = '0.998'
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("D9").Value = '0.710'
$ws.Range("E9").Value = '  -4.54%  '
$ws.Range("E10").Value = '  -6.41%  '
$ws.Range("D11").Value = '0.0000327'
$ws.Range("E11").Value = '  -8.40%  '
$ws.Range("D12").Value = '41.51'
$ws.Range("E12").Value = '  -3.93%  '
$ws.Range("D13").Value = '4.462.75'
$ws.Range("E13").Value = '  -1.27%  '
$ws.Range("D14").Value = '10.10'
$ws.Range("E14").Value = '  -3.84%  '
$ws.Range("D15").Value = '3.857.66'
$ws.Range("E15").Value = '  -1.14%  '
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").Value = '20.64'
$ws.Range("E16").Value = '  +3.14%  '
$ws.Range("B17").Value = 'Uniswap'
$ws.Range("C17").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D17").Value = '13.82'
$ws.Range("E17").Value = '  -2.96%  '
$ws.Range("E18").Value = '  -1.57%  '
$ws.Range("D19").Value = '1.20'
$ws.Range("E19").Value = '  +4.73%  '
$ws.Range("D20").Value = '68.686.65'
$ws.Range("E20").Value = '  +0.40%  '
$ws.Range("D21").Value = '418.77'
$ws.Range("E21").Value = '  -3.39%  '
$ws.Range("E22").Value = '  -5.40%  '
$ws.Range("D23").Value = '14.01'
$ws.Range("E23").Value = '  -4.74%  '
$ws.Range("D24").Value = '86.70'
$ws.Range("E24").Value = '  -3.69%  '
$ws.Range("D25").Value = '3.94'
$ws.Range("E25").Value = '  +5.10%  '
$ws.Range("D26").Value = '11.26'
$ws.Range("E26").Value = '  -8.89%  '
$ws.Range("D27").Value = '10.47'
$ws.Range("E27").Value = '  -4.66%  '
$ws.Range("D28").Value = '35.86'
$ws.Range("E28").Value = '  -4.48%  '
$ws.Range("D29").Value = '682.52'
$ws.Range("E29").Value = '  -4.20%  '
$ws.Range("D30").Value = '13.04'
$ws.Range("E30").Value = '  -2.76%  '
$ws.Range("E31").Value = '  -5.16%  '
$ws.Range("D32").Value = '2.82'
$ws.Range("E32").Value = '  -3.36%  '
$ws.Range("D33").Value = '67.55'
$ws.Range("E33").Value = '  +9.80%  '
$ws.Range("E34").Value = '  +6.33%  '
$ws.Range("D35").Value = '5.86'
$ws.Range("E35").Value = '  -3.65%  '
$ws.Range("D36").Value = '0.0₃0842'
$ws.Range("E36").Value = '  -6.41%  '
$ws.Range("D37").Value = '39.44'
$ws.Range("E37").Value = '  -3.39%  '
$ws.Range("E38").Value = '  -0.06%  '
$ws.Range("E39").Value = '  -1.22%  '
$ws.Range("E40").Value = '  -0.07%  '
$ws.Range("D41").Value = '3.19'
$ws.Range("E41").Value = '  +3.04%  '
$ws.Range("D42").Value = '0.0476'
$ws.Range("E42").Value = '  -3.84%  '
$ws.Range("E43").Value = '  +2.93%  '
$ws.Range("D44").Value = '2.73'
$ws.Range("E44").Value = '  -8.12%  '
$ws.Range("D45").Value = '3.38'
$ws.Range("E45").Value = '  -0.03%  '
$ws.Range("E46").Value = '  -3.17%  '
$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D47").Value = '2.756.79'
$ws.Range("E47").Value = '  +14.09%  '
$ws.Range("B48").Value = 'Stacks'
$ws.Range("C48").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D48").Value = '2.93'
$ws.Range("E48").Value = '  +3.99%  '
$ws.Range("D49").Value = '0.000269'
$ws.Range("E49").Value = '  +9.45%  '
$ws.Range("D50").Value = '143.55'
$ws.Range("E50").Value = '  +0.57%  '
$ws.Range("D51").Value = '3.24'
$ws.Range("E51").Value = '  -3.80%  '
